# Fruta / hortaliza, semanal
# Insert two new weekly observation rows at the top of the Achicoria data
# block (rows 24-25), pushing the existing rows down by two. This mirrors
# the author's edit: two brand-new rows are added and everything that used
# to live at row N (24 <= N <= 94) now lives at row N+2 (26 <= N+2 <= 96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 24..94 down to 26..96, making room for the two
# new rows. Excel's Rows.Insert() also carries the row-above formatting
# (including the date number format on column D) onto the freshly inserted
# rows, just like a manual "Insert Copied Cells" / "Insert Sheet Rows" in
# the UI would.
$ws.Rows("24:25").Insert()

# --- New row 24 ---
$ws.Cells.Item(24, 1).Value = 10
$ws.Cells.Item(24, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(24, 3).Value = "La Araucanía"
$ws.Cells.Item(24, 4).Value = 45054
$ws.Cells.Item(24, 5).Value = 9
$ws.Cells.Item(24, 6).Value = 100112010
$ws.Cells.Item(24, 7).Value = "Achicoria"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 200
$ws.Cells.Item(24, 11).Value = 10000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 10000
$ws.Cells.Item(24, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 556
$ws.Cells.Item(24, 17).Value = 18
$ws.Cells.Item(24, 18).Value = "Hortaliza"

# --- New row 25 ---
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 45054
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = 100112010
$ws.Cells.Item(25, 7).Value = "Achicoria"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 8000
$ws.Cells.Item(25, 12).Value = 8000
$ws.Cells.Item(25, 13).Value = 8000
$ws.Cells.Item(25, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 444
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"
